$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 299.9  # H6: was 225.125
$ws.Cells.Item(6, 9).Value = 311  # I6: was 225.125
$ws.Cells.Item(6, 10).Value = 200  # J6: was 0
$ws.Cells.Item(6, 11).Value = 933  # K6: was 675.375
$ws.Cells.Item(6, 12).Value = 600  # L6: was 0
$ws.Cells.Item(6, 13).Value = -821  # M6: was -563.375
$ws.Cells.Item(6, 14).Value = -824  # N6: was None
$ws.Cells.Item(18, 8).Value = 1535.5714  # H18: was 1456
$ws.Cells.Item(18, 9).Value = 1535.5714  # I18: was 1456
$ws.Cells.Item(18, 11).Value = 1535.5714  # K18: was 1456
$ws.Cells.Item(18, 13).Value = -1251.5714  # M18: was -1172
$ws.Cells.Item(62, 8).Value = 2849.2856  # H62: was 2874.5
$ws.Cells.Item(62, 9).Value = 2849.3333  # I62: was 2879.6
$ws.Cells.Item(62, 11).Value = 2849.3333  # K62: was 2879.6
$ws.Cells.Item(62, 13).Value = -2225.3333  # M62: was -2255.6
$ws.Cells.Item(65, 8).Value = 2849.2856  # H65: was 2874.5
$ws.Cells.Item(65, 9).Value = 2849.3333  # I65: was 2879.6
$ws.Cells.Item(65, 11).Value = 14246.6665  # K65: was 14398
$ws.Cells.Item(65, 13).Value = -11126.6665  # M65: was -11278
$ws.Cells.Item(100, 8).Value = 945.44446  # H100: was 1223.1111
$ws.Cells.Item(100, 9).Value = 387  # I100: was 401.66666
$ws.Cells.Item(100, 10).Value = 2900  # J100: was 2866
$ws.Cells.Item(100, 11).Value = 387  # K100: was 401.66666
$ws.Cells.Item(100, 12).Value = 2900  # L100: was 2866
$ws.Cells.Item(100, 13).Value = 154  # M100: was 139.33334
$ws.Cells.Item(100, 14).Value = -3982  # N100: was -3948
$ws.Cells.Item(111, 8).Value = 3426.4285  # H111: was 3426.5715
$ws.Cells.Item(111, 10).Value = 5658  # J111: was 5658.3335
$ws.Cells.Item(111, 12).Value = 16974  # L111: was 16975.0005
$ws.Cells.Item(111, 14).Value = -23108  # N111: was -23109.0005
$ws.Cells.Item(138, 8).Value = 5143.3213  # H138: was 5075.4546
$ws.Cells.Item(138, 10).Value = 5564.9473  # J138: was 5644.6924
$ws.Cells.Item(138, 12).Value = 16694.8419  # L138: was 16934.0772
$ws.Cells.Item(138, 14).Value = -26974.8419  # N138: was -27214.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 6207.4  # H45: was 8276.714
$ws.Cells.Item(45, 9).Value = 6207.4  # I45: was 8276.714
$ws.Cells.Item(45, 11).Value = 6207.4  # K45: was 8276.714
$ws.Cells.Item(45, 13).Value = -5830.4  # M45: was -7899.714
$ws.Cells.Item(61, 8).Value = 1518.875  # H61: was 1744.5
$ws.Cells.Item(61, 9).Value = 1521.7142  # I61: was 1738.4
$ws.Cells.Item(61, 10).Value = 1499  # J61: was 1775
$ws.Cells.Item(61, 11).Value = 1521.7142  # K61: was 1738.4
$ws.Cells.Item(61, 12).Value = 1499  # L61: was 1775
$ws.Cells.Item(61, 13).Value = -1309.7142  # M61: was -1526.4
$ws.Cells.Item(61, 14).Value = -1923  # N61: was -2199
$ws.Cells.Item(74, 8).Value = 987.5  # H74: was 940.3077
$ws.Cells.Item(74, 9).Value = 893.8570999999999  # I74: was 860.6
$ws.Cells.Item(74, 11).Value = 893.8570999999999  # K74: was 860.6
$ws.Cells.Item(74, 13).Value = -19.85709999999995  # M74: was 13.39999999999998
$ws.Cells.Item(77, 8).Value = 987.5  # H77: was 940.3077
$ws.Cells.Item(77, 9).Value = 893.8570999999999  # I77: was 860.6
$ws.Cells.Item(77, 11).Value = 4469.2855  # K77: was 4303
$ws.Cells.Item(77, 13).Value = -101.2855  # M77: was 65
$ws.Cells.Item(97, 8).Value = 813  # H97: was 1211.625
$ws.Cells.Item(97, 9).Value = 590.6667  # I97: was 1323.5
$ws.Cells.Item(97, 10).Value = 1079.8  # J97: was 1099.75
$ws.Cells.Item(97, 11).Value = 590.6667  # K97: was 1323.5
$ws.Cells.Item(97, 12).Value = 1079.8  # L97: was 1099.75
$ws.Cells.Item(97, 13).Value = -94.66669999999999  # M97: was -827.5
$ws.Cells.Item(97, 14).Value = -2071.8  # N97: was -2091.75
$ws.Cells.Item(132, 8).Value = 1363  # H132: was 1731
$ws.Cells.Item(132, 9).Value = 831.6667  # I132: was 995
$ws.Cells.Item(132, 10).Value = 1761.5  # J132: was 2099
$ws.Cells.Item(132, 11).Value = 2495.0001  # K132: was 2985
$ws.Cells.Item(132, 12).Value = 5284.5  # L132: was 6297
$ws.Cells.Item(132, 13).Value = 34.9998999999998  # M132: was -455
$ws.Cells.Item(132, 14).Value = -10344.5  # N132: was -11357
$ws.Cells.Item(136, 8).Value = 1518.875  # H136: was 1744.5
$ws.Cells.Item(136, 9).Value = 1521.7142  # I136: was 1738.4
$ws.Cells.Item(136, 10).Value = 1499  # J136: was 1775
$ws.Cells.Item(136, 11).Value = 4565.142599999999  # K136: was 5215.200000000001
$ws.Cells.Item(136, 12).Value = 4497  # L136: was 5325
$ws.Cells.Item(136, 13).Value = -2015.142599999999  # M136: was -2665.200000000001
$ws.Cells.Item(136, 14).Value = -9597  # N136: was -10425

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7182.5  # H20: was 7316.8887
$ws.Cells.Item(20, 9).Value = 6839.5713  # I20: was 6984
$ws.Cells.Item(20, 11).Value = 6839.5713  # K20: was 6984
$ws.Cells.Item(20, 13).Value = -6592.5713  # M20: was -6737
$ws.Cells.Item(94, 8).Value = 1123.7778  # H94: was 1250.625
$ws.Cells.Item(94, 9).Value = 587.8570999999999  # I94: was 667.6667
$ws.Cells.Item(94, 11).Value = 587.8570999999999  # K94: was 667.6667
$ws.Cells.Item(94, 13).Value = -136.8570999999999  # M94: was -216.6667
$ws.Cells.Item(107, 8).Value = 6263.5713  # H107: was 7141.857
$ws.Cells.Item(107, 9).Value = 5569.2  # I107: was 6798.8
$ws.Cells.Item(107, 11).Value = 5569.2  # K107: was 6798.8
$ws.Cells.Item(107, 13).Value = -3649.2  # M107: was -4878.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2672.923  # H22: was 2818.0908
$ws.Cells.Item(22, 9).Value = 2879.9  # I22: was 2880
$ws.Cells.Item(22, 10).Value = 1983  # J22: was 2199
$ws.Cells.Item(22, 11).Value = 2879.9  # K22: was 2880
$ws.Cells.Item(22, 12).Value = 1983  # L22: was 2199
$ws.Cells.Item(22, 13).Value = -2529.9  # M22: was -2530
$ws.Cells.Item(22, 14).Value = -2683  # N22: was -2899
$ws.Cells.Item(31, 8).Value = 2104.7083  # H31: was 2135.04
$ws.Cells.Item(31, 9).Value = 1472.8  # I31: was 1525.4445
$ws.Cells.Item(31, 10).Value = 2556.0715  # J31: was 2477.9375
$ws.Cells.Item(31, 11).Value = 1472.8  # K31: was 1525.4445
$ws.Cells.Item(31, 12).Value = 2556.0715  # L31: was 2477.9375
$ws.Cells.Item(31, 13).Value = -1177.8  # M31: was -1230.4445
$ws.Cells.Item(31, 14).Value = -3146.0715  # N31: was -3067.9375
$ws.Cells.Item(34, 8).Value = 2104.7083  # H34: was 2135.04
$ws.Cells.Item(34, 9).Value = 1472.8  # I34: was 1525.4445
$ws.Cells.Item(34, 10).Value = 2556.0715  # J34: was 2477.9375
$ws.Cells.Item(34, 11).Value = 1472.8  # K34: was 1525.4445
$ws.Cells.Item(34, 12).Value = 2556.0715  # L34: was 2477.9375
$ws.Cells.Item(34, 13).Value = -1270.8  # M34: was -1323.4445
$ws.Cells.Item(34, 14).Value = -2960.0715  # N34: was -2881.9375
$ws.Cells.Item(132, 8).Value = 4308.75  # H132: was 4940
$ws.Cells.Item(132, 9).Value = 4782  # I132: was 5433
$ws.Cells.Item(132, 11).Value = 14346  # K132: was 16299
$ws.Cells.Item(132, 13).Value = -11816  # M132: was -13769
$ws.Cells.Item(134, 8).Value = 2607.158  # H134: was 2632.2354
$ws.Cells.Item(134, 9).Value = 2618.6667  # I134: was 2632.2354
$ws.Cells.Item(134, 10).Value = 2400  # J134: was 0
$ws.Cells.Item(134, 11).Value = 7856.000100000001  # K134: was 7896.706200000001
$ws.Cells.Item(134, 12).Value = 7200  # L134: was 0
$ws.Cells.Item(134, 13).Value = -5321.000100000001  # M134: was -5361.706200000001
$ws.Cells.Item(134, 14).Value = -12270  # N134: was None

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 2707  # H38: was 3141.5
$ws.Cells.Item(38, 10).Value = 318.75  # J38: was 391.66666
$ws.Cells.Item(38, 12).Value = 956.25  # L38: was 1174.99998
$ws.Cells.Item(38, 14).Value = -1650.25  # N38: was -1868.99998
$ws.Cells.Item(40, 8).Value = 105  # H40: was 130.75
$ws.Cells.Item(40, 9).Value = 61.11111  # I40: was 57
$ws.Cells.Item(40, 10).Value = 500  # J40: was 499.5
$ws.Cells.Item(40, 11).Value = 244.44444  # K40: was 228
$ws.Cells.Item(40, 12).Value = 2000  # L40: was 1998
$ws.Cells.Item(40, 13).Value = -175.44444  # M40: was -159
$ws.Cells.Item(40, 14).Value = -2138  # N40: was -2136
$ws.Cells.Item(68, 8).Value = 1886.0869  # H68: was 1854.32
$ws.Cells.Item(68, 10).Value = 1903.6818  # J68: was 1869.125
$ws.Cells.Item(68, 12).Value = 5711.0454  # L68: was 5607.375
$ws.Cells.Item(68, 14).Value = -7333.0454  # N68: was -7229.375
$ws.Cells.Item(71, 8).Value = 1886.0869  # H71: was 1854.32
$ws.Cells.Item(71, 10).Value = 1903.6818  # J71: was 1869.125
$ws.Cells.Item(71, 12).Value = 17133.1362  # L71: was 16822.125
$ws.Cells.Item(71, 14).Value = -25245.1362  # N71: was -24934.125
$ws.Cells.Item(88, 8).Value = 20000  # H88: was 18999.8
$ws.Cells.Item(88, 10).Value = 20000  # J88: was 18999.8
$ws.Cells.Item(88, 12).Value = 60000  # L88: was 56999.39999999999
$ws.Cells.Item(88, 14).Value = -60856  # N88: was -57855.39999999999
$ws.Cells.Item(91, 8).Value = 20000  # H91: was 18999.8
$ws.Cells.Item(91, 10).Value = 20000  # J91: was 18999.8
$ws.Cells.Item(91, 12).Value = 60000  # L91: was 56999.39999999999
$ws.Cells.Item(91, 14).Value = -62964  # N91: was -59963.39999999999
$ws.Cells.Item(98, 8).Value = 2728.7778  # H98: was 2882.375
$ws.Cells.Item(98, 9).Value = 2547.25  # I98: was 2896.3333
$ws.Cells.Item(98, 11).Value = 7641.75  # K98: was 8688.999899999999
$ws.Cells.Item(98, 13).Value = -6143.75  # M98: was -7190.999899999999
$ws.Cells.Item(112, 8).Value = 15805  # H112: was 20000
$ws.Cells.Item(112, 9).Value = 9512.5  # I112: was 0
$ws.Cells.Item(112, 11).Value = 28537.5  # K112: was 0
$ws.Cells.Item(112, 13).Value = -27429.5  # M112: was None
$ws.Cells.Item(113, 8).Value = 788.5238000000001  # H113: was 851
$ws.Cells.Item(113, 10).Value = 899.17645  # J113: was 993.06665
$ws.Cells.Item(113, 12).Value = 2697.52935  # L113: was 2979.19995
$ws.Cells.Item(113, 14).Value = -7037.529350000001  # N113: was -7319.19995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 6962.222  # H126: was 7089
$ws.Cells.Item(126, 9).Value = 2699.6667  # I126: was 3999.5
$ws.Cells.Item(126, 10).Value = 9093.5  # J126: was 8633.75
$ws.Cells.Item(126, 11).Value = 8099.000100000001  # K126: was 11998.5
$ws.Cells.Item(126, 12).Value = 27280.5  # L126: was 25901.25
$ws.Cells.Item(126, 13).Value = -5629.000100000001  # M126: was -9528.5
$ws.Cells.Item(126, 14).Value = -32220.5  # N126: was -30841.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3442  # H132: was 3611.7778
$ws.Cells.Item(132, 9).Value = 3185.8572  # I132: was 3375.75
$ws.Cells.Item(132, 11).Value = 9557.571599999999  # K132: was 10127.25
$ws.Cells.Item(132, 13).Value = -7027.571599999999  # M132: was -7597.25
$ws.Cells.Item(136, 8).Value = 2144.2727  # H136: was 2370.889
$ws.Cells.Item(136, 9).Value = 1633.3334  # I136: was 1887.75
$ws.Cells.Item(136, 11).Value = 4900.0002  # K136: was 5663.25
$ws.Cells.Item(136, 13).Value = -2350.0002  # M136: was -3113.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 1000  # H22: was 0
$ws.Cells.Item(22, 9).Value = 1000  # I22: was 0
$ws.Cells.Item(22, 11).Value = 1000  # K22: was 0
$ws.Cells.Item(22, 13).Value = -707  # M22: was None
$ws.Cells.Item(68, 8).Value = 0  # H68: was 21249.25
$ws.Cells.Item(68, 9).Value = 0  # I68: was 21666.666
$ws.Cells.Item(68, 10).Value = 0  # J68: was 19997
$ws.Cells.Item(68, 11).Value = 0  # K68: was 21666.666
$ws.Cells.Item(68, 12).Value = 0  # L68: was 19997
$ws.Cells.Item(68, 13).ClearContents()  # M68: was -20855.666
$ws.Cells.Item(68, 14).ClearContents()  # N68: was -21619
$ws.Cells.Item(71, 8).Value = 0  # H71: was 21249.25
$ws.Cells.Item(71, 9).Value = 0  # I71: was 21666.666
$ws.Cells.Item(71, 10).Value = 0  # J71: was 19997
$ws.Cells.Item(71, 11).Value = 0  # K71: was 64999.99800000001
$ws.Cells.Item(71, 12).Value = 0  # L71: was 59991
$ws.Cells.Item(71, 13).ClearContents()  # M71: was -60943.99800000001
$ws.Cells.Item(71, 14).ClearContents()  # N71: was -68103
$ws.Cells.Item(122, 8).Value = 4345.357  # H122: was 5256.25
$ws.Cells.Item(122, 9).Value = 2926.5  # I122: was 3938.125
$ws.Cells.Item(122, 11).Value = 8779.5  # K122: was 11814.375
$ws.Cells.Item(122, 13).Value = -6329.5  # M122: was -9364.375

